# SummaryStats update: only output maxes that are detected.
# Updates "Max Cr" / "Max Date" (sheet "Alluvial for Mapping", cols V/W)
# and the mirrored "Max Cr" / "Max Date" exhibit table
# (sheet "Alluvial Exhibit", cols G/H) so the reported max values reflect
# the true detected maximum (and its sample date) instead of the prior
# placeholder ceiling value of 10.0 / 10.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Alluvial for Mapping")
$ws2 = $wb.Worksheets.Item("Alluvial Exhibit")

# Writes $value into $ws.Range($targetAddr) as literal TEXT, leaving the
# target cell's existing number format / style untouched.
#
# A plain `Range.Value = "8.15"` (or a date-shaped string) gets
# auto-coerced to a Number/Date by COM, which is wrong here - the sheet
# stores these as plain text. To avoid that we stage the text in a
# scratch cell that has been explicitly marked Text ("@") so it is
# stored verbatim, copy just that resulting value onto the real target
# (xlPasteValues - formatting of the target is not touched), then wipe
# the scratch cell.
function Set-TextValue {
    param($ws, [string]$targetAddr, [string]$value)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---- "Alluvial for Mapping" sheet: Max Cr (V) / Max Date (W) ----

Set-TextValue $ws1 "V5"  "8.15"
Set-TextValue $ws1 "W5"  "2013-03-28"

Set-TextValue $ws1 "V6"  "9.43"
Set-TextValue $ws1 "W6"  "2013-03-29"

Set-TextValue $ws1 "V11" "5.76"
Set-TextValue $ws1 "W11" "2010-07-06"

Set-TextValue $ws1 "V12" "5.44"
Set-TextValue $ws1 "W12" "2010-07-07"

Set-TextValue $ws1 "V14" "2.54"
Set-TextValue $ws1 "W14" "2001-11-13"

Set-TextValue $ws1 "V15" "No Detect Data"
Set-TextValue $ws1 "W15" "No Detect Data"

Set-TextValue $ws1 "V17" "3.01"
Set-TextValue $ws1 "W17" "2010-06-08"

# ---- "Alluvial Exhibit" sheet: Max Cr (G) / Max Date (H) ----

Set-TextValue $ws2 "G7"  "8.15"
Set-TextValue $ws2 "H7"  "3/28/13"

Set-TextValue $ws2 "G8"  "9.43"
Set-TextValue $ws2 "H8"  "3/29/13"

Set-TextValue $ws2 "G14" "5.76"
Set-TextValue $ws2 "H14" "7/6/10"

Set-TextValue $ws2 "G15" "5.44"
Set-TextValue $ws2 "H15" "7/7/10"

Set-TextValue $ws2 "G18" "2.54"
Set-TextValue $ws2 "H18" "11/13/01"

Set-TextValue $ws2 "G19" "NA"
Set-TextValue $ws2 "H19" "NA"

Set-TextValue $ws2 "G21" "3.01"
Set-TextValue $ws2 "H21" "6/8/10"
